$wb = $excel.ActiveWorkbook

$ov   = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shown in Overview!E2:F3 and in the "Status" column (C) of
#    both the zh-cn and de-de detail sheets.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ov.Range("E2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn handback finished at 2016-10-20 09:23:20 (K2/K3, same timestamp)
# ---------------------------------------------------------------------------
$zhcn.Range("K2").Value = "2016-10-20 09:23:20"
$zhcn.Range("K3").Value = "2016-10-20 09:23:20"

# ---------------------------------------------------------------------------
# 3. de-de handback finished later, at 2016-10-20 09:23:38 (K2/K3)
# ---------------------------------------------------------------------------
$dede.Range("K2").Value = "2016-10-20 09:23:38"
$dede.Range("K3").Value = "2016-10-20 09:23:38"

# ---------------------------------------------------------------------------
# 4. Latest Handback File (column J) now has the generated xliff file name
#    for both languages (plain text, not a hyperlink).
# ---------------------------------------------------------------------------
$zhcn.Range("J2").Value = "4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.3f3cb1b5b3f25cfda0a328f0a04e8bb0ba0f1b10.zh-cn.xlf"
$zhcn.Range("J3").Value = "71edd57b-94b9-4a04-bed9-e3376c25f4ea.9005a1e8becfcdfddd69fe73c5bc3d0449983383.zh-cn.xlf"
$dede.Range("J2").Value = "4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.3f3cb1b5b3f25cfda0a328f0a04e8bb0ba0f1b10.de-de.xlf"
$dede.Range("J3").Value = "71edd57b-94b9-4a04-bed9-e3376c25f4ea.9005a1e8becfcdfddd69fe73c5bc3d0449983383.de-de.xlf"

# ---------------------------------------------------------------------------
# 5. Latest Target File (column I) now links back to the source .md file,
#    exactly like column A already does. Rebuild the hyperlink collection so
#    the final order is A2, I2, A3, I3 (matches how Excel lists them when the
#    new ones are added after the matching source-file row).
# ---------------------------------------------------------------------------
foreach ($ws in @($zhcn, $dede)) {
    # Engine quirk: Range.Hyperlinks.Delete() clears every hyperlink on the
    # sheet, regardless of which range it was invoked on - use that to reset
    # and then rebuild in the desired order.
    $ws.Range("A1").Hyperlinks.Delete()

    $hls = $ws.Hyperlinks
    $hls.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/069821608d8497df7353be014bf5eca787b3e92f/e2e/4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.md", [System.Type]::Missing, [System.Type]::Missing, "4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.md") | Out-Null
    $hls.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/069821608d8497df7353be014bf5eca787b3e92f/e2e/4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.md", [System.Type]::Missing, [System.Type]::Missing, "4fd95e1b-8e22-49e1-bdba-586ff34a3e9a.md") | Out-Null
    $hls.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/069821608d8497df7353be014bf5eca787b3e92f/e2e/71edd57b-94b9-4a04-bed9-e3376c25f4ea.md", [System.Type]::Missing, [System.Type]::Missing, "71edd57b-94b9-4a04-bed9-e3376c25f4ea.md") | Out-Null
    $hls.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/069821608d8497df7353be014bf5eca787b3e92f/e2e/71edd57b-94b9-4a04-bed9-e3376c25f4ea.md", [System.Type]::Missing, [System.Type]::Missing, "71edd57b-94b9-4a04-bed9-e3376c25f4ea.md") | Out-Null
}

# ---------------------------------------------------------------------------
# 6. Column widths widened to fit the newly generated, longer content.
#    (ColumnWidth is stored in whole "characters" internally and is rounded
#    to the nearest 1/6th of a character by the engine, so the closest
#    achievable value is used.)
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.1666666
$ov.Columns.Item(6).ColumnWidth = 29.1666666

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666
$zhcn.Columns.Item(9).ColumnWidth = 39.1666666
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666

$dede.Columns.Item(3).ColumnWidth = 29.1666666
$dede.Columns.Item(9).ColumnWidth = 39.1666666
$dede.Columns.Item(10).ColumnWidth = 39.1666666
